# Insert a brand new data row at row 35 of the single worksheet (the
# weekly "Apio" / Feria Lagunitas de Puerto Montt price listing), pushing
# the existing rows 35-186 down to 36-187 (dimension grows to A1:R187),
# then populate the newly inserted row with its own record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 35..186 down by one, inserting a fresh blank row 35.
$ws.Rows("35:35").Insert()

# Fill in the new row 35 with the added record.
$ws.Range("A35").Value = 4
$ws.Range("B35").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C35").Value = 'Los Lagos'
$ws.Range("D35").Value = 44575
$ws.Range("E35").Value = 10
$ws.Range("F35").Value = 100112017
$ws.Range("G35").Value = 'Apio'
$ws.Range("H35").Value = 'Americana (o)'
$ws.Range("I35").Value = 'Primera'
$ws.Range("J35").Value = 40
$ws.Range("K35").Value = 12000
$ws.Range("L35").Value = 12000
$ws.Range("M35").Value = 12000
$ws.Range("N35").Value = '$/docena de matas'
$ws.Range("O35").Value = 'Región de Coquimbo'
$ws.Range("P35").Value = 2000
$ws.Range("Q35").Value = 6
$ws.Range("R35").Value = 'Hortaliza'
